# Update "想去人数" (interest count, column F) figures for several rows
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types)
# sheets, matching the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet "展览") ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 97
$ws1.Range("F8").Value  = 130
$ws1.Range("F13").Value = 156
$ws1.Range("F14").Value = 1375
$ws1.Range("F16").Value = 199
$ws1.Range("F17").Value = 312
$ws1.Range("F19").Value = 717
$ws1.Range("F23").Value = 2553
$ws1.Range("F24").Value = 1312
$ws1.Range("F26").Value = 224
$ws1.Range("F28").Value = 959
$ws1.Range("F30").Value = 1104
$ws1.Range("F33").Value = 759
$ws1.Range("F34").Value = 449

# --- 演出 (sheet "演出") ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 177
$ws2.Range("F13").Value = 534

# --- 全部类型 (sheet "全部类型") ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 177
$ws4.Range("F12").Value = 97
$ws4.Range("F13").Value = 130
$ws4.Range("F20").Value = 156
$ws4.Range("F21").Value = 1375
$ws4.Range("F23").Value = 199
$ws4.Range("F24").Value = 312
$ws4.Range("F27").Value = 2553
$ws4.Range("F29").Value = 1312
$ws4.Range("F34").Value = 224
$ws4.Range("F36").Value = 959
$ws4.Range("F40").Value = 1104
$ws4.Range("F41").Value = 759
$ws4.Range("F42").Value = 449
